$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.946.21'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '1.817.72'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('E4').Value = '  -0.70%  '
$ws.Range('D5').Value = "'311.57"
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').Value = "'0.4294"
$ws.Range('E7').Value = '  +1.42%  '
$ws.Range('E8').Value = '  +2.20%  '
$ws.Range('D9').Value = "'0.07247"
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('E10').Value = '  +2.22%  '
$ws.Range('D11').Value = '2.023.61'
$ws.Range('E11').Value = '  +11.23%  '
$ws.Range('E12').Value = '  +3.78%  '
$ws.Range('D13').Value = "'6.651"
$ws.Range('E13').Value = '  +4.04%  '
$ws.Range('D14').Value = "'5.388"
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').Value = "'0.06918"
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').Value = "'80.64"
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = "'1.005"
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = "'0.000008920"
$ws.Range('E18').Value = '  +2.38%  '
$ws.Range('D19').Value = "'1.005"
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').Value = "'15.26"
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').Value = '26.995.49'
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').Value = "'5.185"
$ws.Range('E22').Value = '  +2.29%  '
$ws.Range('D23').Value = "'11.10"
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '2.253.28'
$ws.Range('E24').Value = '  +10.35%  '
$ws.Range('D25').Value = "'153.89"
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = "'1.883"
$ws.Range('E26').Value = '  -4.01%  '
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').Value = "'5.217"
$ws.Range('E28').Value = '  +3.41%  '
$ws.Range('D29').Value = "'1.887"
$ws.Range('E29').Value = '  +16.72%  '
$ws.Range('D30').Value = "'115.12"
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').Value = "'0.08961"
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('D32').Value = "'1.166"
$ws.Range('E32').Value = '  +6.69%  '
$ws.Range('D33').Value = "'0.7426"
$ws.Range('E33').Value = '  +2.48%  '
$ws.Range('D34').Value = "'4.417"
$ws.Range('E34').Value = '  +1.88%  '
$ws.Range('D35').Value = "'2.802"
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('D37').Value = "'1.125"
$ws.Range('E37').Value = '  +3.48%  '
$ws.Range('D38').Value = "'0.05219"
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('E39').Value = '  +1.21%  '
$ws.Range('D40').Value = "'0.5076"
$ws.Range('E40').Value = '  +2.05%  '
$ws.Range('D41').Value = "'2.731"
$ws.Range('E41').Value = '  +7.93%  '
$ws.Range('D42').Value = "'0.1642"
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('D43').Value = "'6.441"
$ws.Range('E43').Value = '  +7.48%  '
$ws.Range('D44').Value = "'8.233"
$ws.Range('E44').Value = '  +3.97%  '
$ws.Range('D45').Value = "'107.13"
$ws.Range('E45').Value = '  +2.13%  '
$ws.Range('D46').Value = "'10.37"
$ws.Range('E46').Value = '  +2.46%  '
$ws.Range('D47').Value = "'1.005"
$ws.Range('D48').Value = "'1.658"
$ws.Range('E48').Value = '  +4.82%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value = "'0.4580"
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.06301"
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = "'1.803"
$ws.Range('E51').Value = '  +4.98%  '

# Reset style on text-forced numeric-looking cells so they keep the default
# (un-styled) cell format instead of picking up a quote-prefix style.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
